$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string into a cell as a plain shared-string value
# (never inferred as a number/date/bool, and never stamped with a bonus
# style like quotePrefix) by stuffing it through a text formula and then
# collapsing the formula down to its static value via Copy/PasteSpecial
# (xlPasteValues = -4163).
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Start clean: wipe every value on the sheet so the shared-string table is
# rebuilt from scratch in exactly the order we write things below (existing
# strings are de-duplicated by content, not by position, so the only way to
# control final ordering is to clear first).
$ws.UsedRange.ClearContents() | Out-Null

# ---- header row (row 1) --------------------------------------------------
# Written B..I then A last, matching the target shared-string order.
Set-TextValue $ws.Range("B1") "model"
Set-TextValue $ws.Range("C1") "kfun"
Set-TextValue $ws.Range("D1") "lambda"
Set-TextValue $ws.Range("E1") "kparam"
Set-TextValue $ws.Range("F1") "lr"
Set-TextValue $ws.Range("G1") "index"
Set-TextValue $ws.Range("H1") "train"
Set-TextValue $ws.Range("I1") "validation"
Set-TextValue $ws.Range("A1") "function"

# Header cells all carry the bold/bordered/centered header style (style
# index 1 in the original workbook) - copy it from A1 (which already has
# it) across the rest of row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:I1").PasteSpecial(-4122) | Out-Null

# ---- data rows (rows 2-6), written column by column (A, then B, ... I) ---
Set-TextValue $ws.Range("A2") "color_histogram"
Set-TextValue $ws.Range("A3") "cooccurrence_matrix"
Set-TextValue $ws.Range("A4") "deepfeatures"
Set-TextValue $ws.Range("A5") "edge_direction_histogram"
Set-TextValue $ws.Range("A6") "rgb_cooccurrence_matrix"

Set-TextValue $ws.Range("B2") "ksvm"
Set-TextValue $ws.Range("B3") "ksvm"
Set-TextValue $ws.Range("B4") "ksvm"
Set-TextValue $ws.Range("B5") "ksvm"
Set-TextValue $ws.Range("B6") "ksvm"

Set-TextValue $ws.Range("C2") "rbf"
Set-TextValue $ws.Range("C3") "rbf"
Set-TextValue $ws.Range("C4") "rbf"
Set-TextValue $ws.Range("C5") "rbf"
Set-TextValue $ws.Range("C6") "rbf"

Set-TextValue $ws.Range("D2") "0.3"
Set-TextValue $ws.Range("D3") "0.1"
Set-TextValue $ws.Range("D4") "0.3"
Set-TextValue $ws.Range("D5") "0.002"
Set-TextValue $ws.Range("D6") "0.003"

Set-TextValue $ws.Range("E2") "0.1"
Set-TextValue $ws.Range("E3") "0.2"
Set-TextValue $ws.Range("E4") "0.3"
Set-TextValue $ws.Range("E5") "0.2"
Set-TextValue $ws.Range("E6") "0.1"

$ws.Range("F2").Value2 = 0.1
$ws.Range("F3").Value2 = 0.1
$ws.Range("F4").Value2 = 0.1
$ws.Range("F5").Value2 = 0.1
$ws.Range("F6").Value2 = 0.1

Set-TextValue $ws.Range("G2") "ksvm_kfun=rbf_lambda_=0.3_kparam=0.001_lr=0.01__color_histogram"
Set-TextValue $ws.Range("G3") "ksvm_kfun=rbf_lambda_=0.1_kparam=0.003_lr=0.003__cooccurrence_matrix"
Set-TextValue $ws.Range("G4") "ksvm_kfun=rbf_lambda_=0.3_kparam=0.001_lr=0.01__deepfeatures"
Set-TextValue $ws.Range("G5") "ksvm_kfun=rbf_lambda_=0.002_kparam=0.2_lr=0.03__edge_direction_histogram"
Set-TextValue $ws.Range("G6") "ksvm_kfun=rbf_lambda_=0.003_kparam=0.1_lr=0.003__rgb_cooccurrence_matrix"

$ws.Range("H2").Value2 = 0.344279661016949
$ws.Range("H3").Value2 = 0.344279661016949
$ws.Range("H4").Value2 = 0.9242584745762711
$ws.Range("H5").Value2 = 0.344279661016949
$ws.Range("H6").Value2 = 0.344279661016949

$ws.Range("I2").Value2 = 0.333333333333333
$ws.Range("I3").Value2 = 0.333333333333333
$ws.Range("I4").Value2 = 0.9666666666666661
$ws.Range("I5").Value2 = 0.333333333333333
$ws.Range("I6").Value2 = 0.333333333333333
